{"js": "// Artigo 3 (Fa\u00e7ade) \u2013 update the subtitle after \"Nhibernate \u2013 \" from\n// \"introdu\u00e7\u00e3o pr\u00e1tica aos Relacionamentos 1-1 e 1-N\" to\n// \"Relacionamentos 1-1, 1-N, Persist\u00eancia e Padr\u00f5es de Projeto\".\nconst oldSuffix = \" \u2013 introdu\u00e7\u00e3o pr\u00e1tica aos Relacionamentos 1-1 e 1-N\";\nconst newSuffix = \" \u2013 Relacionamentos 1-1, 1-N, Persist\u00eancia e Padr\u00f5es de Projeto\";\n\nconst body = context.document.body;\nconst results = body.search(oldSuffix, { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the title suffix to update.\");\n}\n\n// Replace in place so the run keeps its existing character formatting\n// (Verdana font) instead of inheriting formatting from wherever the cursor\n// would otherwise land.\nresults.items[0].insertText(newSuffix, \"Replace\");\nawait context.sync();\n", "ps1": "# Artigo 3 (Fa\u00e7ade) \u2013 update the subtitle after \"Nhibernate \u2013 \" from\n# \"introdu\u00e7\u00e3o pr\u00e1tica aos Relacionamentos 1-1 e 1-N\" to\n# \"Relacionamentos 1-1, 1-N, Persist\u00eancia e Padr\u00f5es de Projeto\".\n$d = $word.ActiveDocument\n\n$oldSuffix = \" \u2013 introdu\u00e7\u00e3o pr\u00e1tica aos Relacionamentos 1-1 e 1-N\"\n$newSuffix = \" \u2013 Relacionamentos 1-1, 1-N, Persist\u00eancia e Padr\u00f5es de Projeto\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = $oldSuffix\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newSuffix\n$find.Forward = $true\n$find.Wrap = 0\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
